# Master_table.xlsx - "added dish weights for half of the samples 20191029
# to start AFDW process - samples in 75C for 12-24 hrs"
#
# Fills in the wt.dish.DRY.SAMPLE (column N) measurements for rows 2-46.
# Two rows (9 and 16) are marked "-" (no reading for those samples),
# matching the "-" notation already used elsewhere in those rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("N2").Value = 0.45829999999999999
$ws.Range("N3").Value = 0.51529999999999998
$ws.Range("N4").Value = 0.5393
$ws.Range("N5").Value = 0.53539999999999999
$ws.Range("N6").Value = 0.51470000000000005
$ws.Range("N7").Value = 0.49969999999999998
$ws.Range("N8").Value = 0.45429999999999998
$ws.Range("N9").Value = "-"
$ws.Range("N10").Value = 0.46679999999999999
$ws.Range("N11").Value = 0.42430000000000001
$ws.Range("N12").Value = 0.50639999999999996
$ws.Range("N13").Value = 0.47920000000000001
$ws.Range("N14").Value = 0.48499999999999999
$ws.Range("N15").Value = 0.59279999999999999
$ws.Range("N16").Value = "-"
$ws.Range("N17").Value = 0.4662
$ws.Range("N18").Value = 0.35239999999999999
$ws.Range("N19").Value = 0.377
$ws.Range("N20").Value = 0.48120000000000002
$ws.Range("N21").Value = 0.38379999999999997
$ws.Range("N22").Value = 0.36749999999999999
$ws.Range("N23").Value = 0.43440000000000001
$ws.Range("N24").Value = 0.43369999999999997
$ws.Range("N25").Value = 0.55820000000000003
$ws.Range("N26").Value = 0.4718
$ws.Range("N27").Value = 0.55210000000000004
$ws.Range("N28").Value = 0.57709999999999995
$ws.Range("N29").Value = 0.4844
$ws.Range("N30").Value = 0.36370000000000002
$ws.Range("N31").Value = 0.43769999999999998
$ws.Range("N32").Value = 0.49299999999999999
$ws.Range("N33").Value = 0.44819999999999999
$ws.Range("N34").Value = 0.50580000000000003
$ws.Range("N35").Value = 0.45600000000000002
$ws.Range("N36").Value = 0.40239999999999998
$ws.Range("N37").Value = 0.4012
$ws.Range("N38").Value = 0.35420000000000001
$ws.Range("N39").Value = 0.47939999999999999
$ws.Range("N40").Value = 0.43640000000000001
$ws.Range("N41").Value = 0.46489999999999998
$ws.Range("N42").Value = 0.47170000000000001
$ws.Range("N43").Value = 0.43680000000000002
$ws.Range("N44").Value = 0.4612
$ws.Range("N45").Value = 0.50960000000000005
$ws.Range("N46").Value = 0.46460000000000001

# Leave the view/selection the way the author's session ended up: cell
# N47 selected (last cell touched while entering the new readings).
$ws.Activate()
$ws.Range("N47").Select()
